$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 113 ---
$ws.Range("A113").Value = "What is the maximum number of tracks that can be specified within a single ODF file?"
$ws.Range("B113").Value = "llama3.2:latest"
$ws.Range("C113").Value = "The maximum number of tracks that can be specified within a single ODF file is 200."

# --- Row 114 ---
$ws.Range("A114").Value = "What’s the track limit for an ODF?"

$ws.Range("B114").Value = "llama3.2:latest"

$c114text = @'
The maximum number of tracks that can be specified within a single ODF file is not explicitly stated in the provided documentation. However, it is mentioned that there is a limit of 200 tracks at most.
It's recommended to use OIF files instead of sending unnecessary data by saving only the specific depth range of interest. To do this:
1. Enter the From and To index range on the Depth tab (e.g., 7000-8000).
2. Select "Save Displayed Depth or Time interval As (OIF)" from the GEO menu.
3. Check the Statistics dialog box to ensure that no limits have been exceeded.
Note: The exact track limit for an ODF is not specified in the provided documentation, and it's recommended to use OIF files for more efficient data management.
'@
$ws.Range("C114").Value = $c114text

# --- Widen column C to fit the new, longer responses (561 -> 733 stored width units) ---
$ws.Columns.Item(3).ColumnWidth = 732.17
